# Apply updated cryptocurrency price/volume data scraped on Mon Feb 20 03:45:11 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.392.18"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.683.21"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.26"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3880"
$ws.Range("E7").Value = "  -1.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3999"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.479"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.22"
$ws.Range("E11").Value = "  -3.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08738"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.93"
$ws.Range("E13").Value = "  +10.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.464"
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.975"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001339"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.671.85"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.58"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07189"
$ws.Range("E19").Value = "  +2.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.63"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.234"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  -2.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.384.23"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.008"
$ws.Range("E25").Value = "  -7.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.341"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.45"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.53"
$ws.Range("E28").Value = "  +4.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.603"
$ws.Range("E29").Value = "  +10.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.354"
$ws.Range("E30").Value = "  +3.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "137.91"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.853.80"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08736"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.331"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  -2.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02971"
$ws.Range("E36").Value = "  +8.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.970"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2742"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.74"
$ws.Range("E39").Value = "  -4.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09130"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7949"
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.01"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.467"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.26"
$ws.Range("E44").Value = "  +9.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7178"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.582"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.257"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.393"
$ws.Range("E48").Value = "  +6.46%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.27"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08039"
$ws.Range("E51").Value = "  +0.78%  "
